$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price cells that would otherwise be auto-parsed as numbers to remain text
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "65.235.15"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "3.407.09"
$ws.Range("E3").Value = "  -3.45%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "583.48"
$ws.Range("D6").Value = "136.84"
$ws.Range("E6").Value = "  -4.93%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.407.76"
$ws.Range("E8").Value = "  -3.46%  "
$ws.Range("D9").Value = "0.494"
$ws.Range("D10").Value = "7.18"
$ws.Range("E10").Value = "  -7.52%  "
$ws.Range("E11").Value = "  -10.48%  "
$ws.Range("E12").Value = "  -7.88%  "
$ws.Range("D13").Value = "3.985.28"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("E14").Value = "  -11.30%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "0.115"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "65.259.92"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.409.48"
$ws.Range("E17").Value = "  -4.04%  "
$ws.Range("D18").Value = "25.90"
$ws.Range("E18").Value = "  -9.83%  "
$ws.Range("D19").Value = "9.68"
$ws.Range("E19").Value = "  -11.35%  "
$ws.Range("E20").Value = "  -5.63%  "
$ws.Range("D21").Value = "13.51"
$ws.Range("E21").Value = "  -5.91%  "
$ws.Range("D22").Value = "383.34"
$ws.Range("E22").Value = "  -7.59%  "
$ws.Range("E23").Value = "  -7.59%  "
$ws.Range("D24").Value = "72.60"
$ws.Range("E24").Value = "  -6.16%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "3.544.37"
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("D27").Value = "0.0000104"
$ws.Range("E27").Value = "  -10.56%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").Value = "7.03"
$ws.Range("E29").Value = "  -9.97%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  -10.10%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "8.05"
$ws.Range("E31").Value = "  -9.99%  "
$ws.Range("D32").Value = "3.411.42"
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "0.142"
$ws.Range("E34").Value = "  -8.11%  "
$ws.Range("E35").Value = "  -6.60%  "
$ws.Range("D36").Value = "170.88"
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("D37").Value = "6.74"
$ws.Range("E37").Value = "  -10.51%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.45"
$ws.Range("E38").Value = "  -8.65%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "1.12"
$ws.Range("E39").Value = "  -13.20%  "
$ws.Range("D40").Value = "4.70"
$ws.Range("E40").Value = "  -10.72%  "
$ws.Range("D41").Value = "0.0755"
$ws.Range("E41").Value = "  -7.94%  "
$ws.Range("E42").Value = "  -5.45%  "
$ws.Range("D43").Value = "43.47"
$ws.Range("E43").Value = "  -3.99%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -13.82%  "
$ws.Range("E46").Value = "  -11.40%  "
$ws.Range("D47").Value = "1.08"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").Value = "22.13"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("D49").Value = "6.48"
$ws.Range("E49").Value = "  -8.51%  "
$ws.Range("D50").Value = "2.03"
$ws.Range("E50").Value = "  -15.98%  "
$ws.Range("D51").Value = "2.173.45"
$ws.Range("E51").Value = "  -7.30%  "
